$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.709.66'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +2.49%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').Value = '2.202.40'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E3').ClearFormats()
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').Value = "'257.21"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.27%  '
$ws.Range('E5').ClearFormats()
$ws.Range('D6').Value = "'84.65"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +13.27%  '
$ws.Range('E6').ClearFormats()
$ws.Range('D7').Value = "'0.621"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.94%  '
$ws.Range('E7').ClearFormats()
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E8').ClearFormats()
$ws.Range('E9').Value = '  +1.68%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').Value = "'45.41"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +11.13%  '
$ws.Range('E10').ClearFormats()
$ws.Range('D11').Value = "'0.0922"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.65%  '
$ws.Range('E11').ClearFormats()
$ws.Range('D12').Value = "'7.26"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.79%  '
$ws.Range('E12').ClearFormats()
$ws.Range('E13').Value = '  +2.51%  '
$ws.Range('E13').ClearFormats()
$ws.Range('D14').Value = '2.531.64'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('E14').ClearFormats()
$ws.Range('D15').Value = "'14.40"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.00%  '
$ws.Range('E15').ClearFormats()
$ws.Range('D16').Value = '2.194.63'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('E16').ClearFormats()
$ws.Range('D17').Value = "'0.786"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.54%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').Value = '43.656.86'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.54%  '
$ws.Range('E18').ClearFormats()
$ws.Range('E19').Value = '  +1.08%  '
$ws.Range('E19').ClearFormats()
$ws.Range('D20').Value = "'69.94"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.53%  '
$ws.Range('E20').ClearFormats()
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('E21').ClearFormats()
$ws.Range('E22').Value = '  +10.71%  '
$ws.Range('E22').ClearFormats()
$ws.Range('D23').Value = "'231.81"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.60%  '
$ws.Range('E23').ClearFormats()
$ws.Range('D24').Value = "'8.98"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.10%  '
$ws.Range('E24').ClearFormats()
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E25').ClearFormats()
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('B26').ClearFormats()
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C26').ClearFormats()
$ws.Range('D26').Value = "'10.66"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.14%  '
$ws.Range('E26').ClearFormats()
$ws.Range('B27').Value = 'WEMIXToken'
$ws.Range('B27').ClearFormats()
$ws.Range('C27').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('C27').ClearFormats()
$ws.Range('D27').Value = "'3.50"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.80%  '
$ws.Range('E27').ClearFormats()
$ws.Range('E28').Value = '  +3.15%  '
$ws.Range('E28').ClearFormats()
$ws.Range('D29').Value = "'38.83"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.60%  '
$ws.Range('E29').ClearFormats()
$ws.Range('D30').Value = "'2.23"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.13%  '
$ws.Range('E30').ClearFormats()
$ws.Range('D31').Value = "'173.72"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('E31').ClearFormats()
$ws.Range('E32').Value = '  +1.36%  '
$ws.Range('E32').ClearFormats()
$ws.Range('E33').Value = '  +3.28%  '
$ws.Range('E33').ClearFormats()
$ws.Range('D34').Value = "'5.32"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.47%  '
$ws.Range('E34').ClearFormats()
$ws.Range('E35').Value = '  +1.77%  '
$ws.Range('E35').ClearFormats()
$ws.Range('D36').Value = "'0.111"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.78%  '
$ws.Range('E36').ClearFormats()
$ws.Range('D37').Value = "'0.0360"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +5.10%  '
$ws.Range('E37').ClearFormats()
$ws.Range('E38').Value = '  +5.14%  '
$ws.Range('E38').ClearFormats()
$ws.Range('D39').Value = "'12.62"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.85%  '
$ws.Range('E39').ClearFormats()
$ws.Range('D40').Value = "'2.87"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +5.32%  '
$ws.Range('E40').ClearFormats()
$ws.Range('D41').Value = "'2.09"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('E41').ClearFormats()
$ws.Range('D42').Value = "'63.60"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +6.57%  '
$ws.Range('E42').ClearFormats()
$ws.Range('E43').Value = '  +4.38%  '
$ws.Range('E43').ClearFormats()
$ws.Range('E44').Value = '  +1.79%  '
$ws.Range('E44').ClearFormats()
$ws.Range('B45').Value = 'Aave'
$ws.Range('B45').ClearFormats()
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C45').ClearFormats()
$ws.Range('D45').Value = "'100.33"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.04%  '
$ws.Range('E45').ClearFormats()
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('B46').ClearFormats()
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C46').ClearFormats()
$ws.Range('D46').Value = "'8.33"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.13%  '
$ws.Range('E46').ClearFormats()
$ws.Range('D47').Value = "'0.0980"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('E47').ClearFormats()
$ws.Range('E48').Value = '  +4.78%  '
$ws.Range('E48').ClearFormats()
$ws.Range('E49').Value = '  +1.01%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').Value = "'0.434"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -5.19%  '
$ws.Range('E50').ClearFormats()
$ws.Range('E51').Value = '  +7.45%  '
$ws.Range('E51').ClearFormats()
